$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1 holds the document date as a serial number; bump it by one day
# (45310 -> 45311, i.e. 2024-01-19 -> 2024-01-20)
$ws.Range("A1").Value = 45311

# Update the unit price for the three "Soporte LATERAL U" rows
# (185 -> 94.3)
$ws.Range("D27").Value = 94.3
$ws.Range("D28").Value = 94.3
$ws.Range("D29").Value = 94.3
